$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13829
$ws.Range("C3:C4").Value = 12894
$ws.Range("C5:C7").Value = 10907
$ws.Range("C8").Value = 10235
$ws.Range("C9:C13").Value = 9102
$ws.Range("C14:C21").Value = 9077
$ws.Range("C22:C26").Value = 8472
$ws.Range("C27:C28").Value = 8087
$ws.Range("C29:C57").Value = 8042
$ws.Range("C58:C68").Value = 7815
$ws.Range("C69:C93").Value = 7618
$ws.Range("C94:C191").Value = 7534
$ws.Range("C192:C252").Value = 7343
